# Update countries & provincias Spain
# - Refresh the "Nepal" row's stats (it overtakes "Taiwan" in the case-count
#   ranking), and re-sort a handful of near-tied rows (Belice/Santa Lucia,
#   Groenlandia/Montserrat) so the table stays ordered by total cases.
# - Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 132 / 133: Taiwan <-> Nepal -------------------------------------
# Row 132 was Taiwan (440,0,402,31,0,0,7); Row 133 was Nepal (427,0,45,380,0,0,2).
# Nepal's numbers are refreshed and it now ranks above Taiwan, so Nepal takes
# row 132 (with its updated figures) and Taiwan drops to row 133 unchanged.
$ws.Cells.Item(132, 1).Value = "Nepal"
$ws.Cells.Item(132, 2).Value = 444
$ws.Cells.Item(132, 3).Value = 17
$ws.Cells.Item(132, 4).Value = 45
$ws.Cells.Item(132, 5).Value = 397
$ws.Cells.Item(132, 8).Value = 2

$ws.Cells.Item(133, 1).Value = "Taiwan"
$ws.Cells.Item(133, 2).Value = 440
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 402
$ws.Cells.Item(133, 5).Value = 31
$ws.Cells.Item(133, 8).Value = 7

# --- Row 197 / 198: Santa Lucia <-> Belice -------------------------------
$ws.Cells.Item(197, 1).Value = "Belice"
$ws.Cells.Item(197, 4).Value = 16
$ws.Cells.Item(197, 8).Value = 2

$ws.Cells.Item(198, 1).Value = "Santa Lucia"
$ws.Cells.Item(198, 4).Value = 18
$ws.Cells.Item(198, 8).Value = 0

# --- Row 209 / 210: Montserrat <-> Groenlandia ---------------------------
$ws.Cells.Item(209, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 0

$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 8).Value = 1

# --- Timestamp banner -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 05:35"
